$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the target paragraph (Step 3 narrative paragraph that starts
# with "While installation it asks ...") robustly by its text, rather
# than by a hard-coded paragraph index.
# ------------------------------------------------------------------
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "While installation it asks*") {
        $targetIndex = $i
        break
    }
}

$targetPara = $d.Paragraphs($targetIndex)

# ------------------------------------------------------------------
# Step 1: replace the narrative text that precedes the "_GoBack"
# bookmark with the new wording (keeps the bookmark's position intact,
# matching the diff where bookmarkStart/bookmarkEnd sit right after
# the new final run).
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$preRange = $d.Range($targetPara.Range.Start, $bm.Start)
$preRange.Text = "You need to enter the remote machine information in the order shown in the connection-parameters.json file. When you install the Dockerfile it will access the remote machine you entered in the file."

# ------------------------------------------------------------------
# Step 2: the two paragraphs that immediately follow (the "[Unit] ..."
# systemd snippet and the "After=network.target ... [Install] ..."
# continuation) are being removed entirely, along with the obsolete
# tail of the narrative paragraph (" seconds. In the end ... continue.").
# Clear all of that text first (range from the bookmark through the
# end of the second following paragraph).
# ------------------------------------------------------------------
$unitPara = $d.Paragraphs($targetIndex + 1)
$installPara = $d.Paragraphs($targetIndex + 2)
$bm2 = $d.Bookmarks("_GoBack")
$tailRange = $d.Range($bm2.Start, $installPara.Range.End - 1)
$tailRange.Text = ""

# ------------------------------------------------------------------
# Step 3: the paragraph marks can't be removed by deleting a range
# that spans them directly, so collapse each now-empty paragraph into
# its neighbour by deleting just its trailing paragraph mark. Doing
# this from the end backwards means the very last (originally empty,
# unshaded/non-bold) paragraph mark survives and becomes the mark for
# the merged paragraph -- exactly like the target, whose surviving
# paragraph keeps the plain (no shading, no bold) formatting rather
# than the grey-shaded/bold formatting of the deleted service-file
# paragraphs.
# ------------------------------------------------------------------
$installPara = $d.Paragraphs($targetIndex + 2)
$mark = $d.Range($installPara.Range.End - 1, $installPara.Range.End)
$mark.Delete()

$unitPara = $d.Paragraphs($targetIndex + 1)
$mark = $d.Range($unitPara.Range.End - 1, $unitPara.Range.End)
$mark.Delete()

$targetPara = $d.Paragraphs($targetIndex)
$mark = $d.Range($targetPara.Range.End - 1, $targetPara.Range.End)
$mark.Delete()
